# A new weekly price observation row was inserted into the dataset at
# row 660 (1-indexed, matching the worksheet's row numbers). All of the
# existing rows from 660 through 715 shift down by one row (to 661-716),
# and a brand-new row of data is placed at row 660.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 660; this pushes the old
# row 660 (and everything below it) down by one row.
$ws.Rows.Item(660).Insert()

# Populate the newly inserted row 660 with the new observation.
$ws.Cells.Item(660, 1).Value  = 6
$ws.Cells.Item(660, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(660, 3).Value  = "Metropolitana"
$ws.Cells.Item(660, 4).Value  = 45021
$ws.Cells.Item(660, 5).Value  = 13
$ws.Cells.Item(660, 6).Value  = 100112044
$ws.Cells.Item(660, 7).Value  = "Perejil"
$ws.Cells.Item(660, 8).Value  = "Sin especificar"
$ws.Cells.Item(660, 9).Value  = "Primera"
$ws.Cells.Item(660, 10).Value = 290
$ws.Cells.Item(660, 11).Value = 12000
$ws.Cells.Item(660, 12).Value = 13000
$ws.Cells.Item(660, 13).Value = 12448
$ws.Cells.Item(660, 14).Value = "`$/docena de atados"
$ws.Cells.Item(660, 15).Value = "Región Metropolitana"
$ws.Cells.Item(660, 16).Value = 4149
$ws.Cells.Item(660, 17).Value = 3
$ws.Cells.Item(660, 18).Value = "Hortaliza"
